$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.35 = 30242.65 pesos`n✅ 30242.65 pesos = 7.34 = 966.51 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 136
$ws2.Range("O10").Value = 4113
$ws2.Range("N12").Value = 4120.96
$ws2.Range("O12").Value = 131.7
